# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" figures for the first and third data
# rows of the statement (rows 16 and 18) are swapped: what used to be the
# "1812" period with a mora value of 18750 is now the "1810" period with a
# mora value of 3125, and vice versa for the former "1810" row, which now
# shows period "1812" with a mora value of 18750. The middle row (17,
# period "1811") is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 ("Periodo Mora" / "Valor Mora")
$ws.Range("E16").Value = "1810"
$ws.Range("F16").Value = 3125

# Row 18 ("Periodo Mora" / "Valor Mora")
$ws.Range("E18").Value = "1812"
$ws.Range("F18").Value = 18750
